$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four data rows (row 2 through row 5) with corrected values
$ws.Range("B2").Value = "5056780-46.2019.8.21.0001"
$ws.Range("C2").Value = "0140462-81.2019.8.21.0001"
$ws.Range("D2").Value = "CIV.18900.01"

$ws.Range("B3").Value = "5000559-78.2019.8.21.0054"
$ws.Range("C3").Value = "0003337-09.2019.8.21.0054"
$ws.Range("D3").Value = "CIV.08526.01"

$ws.Range("B4").Value = "5001387-66.2023.8.21.9000"
$ws.Range("C4").Value = "5008039-43.2018.8.21.0022"
$ws.Range("D4").Value = "CIV.02502.02"

$ws.Range("B5").Value = "5009986-45.2011.8.21.0001"
$ws.Range("C5").Value = "0323125-76.2011.8.21.0001"
$ws.Range("D5").Value = "CIV.21367.01"

# Remove the now-obsolete rows 6 through 12
$ws.Rows("6:12").Delete()
